# The source data for rows 4, 5 and 6 (the three observation records) was
# re-ordered: the record that used to be on row 5 is now on row 4, the
# record that used to be on row 6 is now on row 5, and the record that
# used to be on row 4 is now on row 6 (a cyclic rotation "up" by one row).
#
# Implement this by reading the full row contents (columns A:AY, which is
# the sheet's used range) for rows 4-6 and writing them back in rotated
# order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng4 = $ws.Range("A4:AY4")
$rng5 = $ws.Range("A5:AY5")
$rng6 = $ws.Range("A6:AY6")

$v4 = $rng4.Value2
$v5 = $rng5.Value2
$v6 = $rng6.Value2

# Columns Y and AA hold free-text dates formatted like "2023-09-04". When
# such text is written back through Value2, Excel happily reinterprets it
# as a real date serial number. Prefix it with an apostrophe (exactly what
# a user would type to force text) so it round-trips as plain text.
$colY = 25   # column Y is the 25th column starting at A
$colAA = 27  # column AA is the 27th column starting at A

if ($v4[1, $colY] -ne $null -and $v4[1, $colY] -ne "") { $v4[1, $colY] = "'" + $v4[1, $colY] }
if ($v4[1, $colAA] -ne $null -and $v4[1, $colAA] -ne "") { $v4[1, $colAA] = "'" + $v4[1, $colAA] }
if ($v5[1, $colY] -ne $null -and $v5[1, $colY] -ne "") { $v5[1, $colY] = "'" + $v5[1, $colY] }
if ($v5[1, $colAA] -ne $null -and $v5[1, $colAA] -ne "") { $v5[1, $colAA] = "'" + $v5[1, $colAA] }
if ($v6[1, $colY] -ne $null -and $v6[1, $colY] -ne "") { $v6[1, $colY] = "'" + $v6[1, $colY] }
if ($v6[1, $colAA] -ne $null -and $v6[1, $colAA] -ne "") { $v6[1, $colAA] = "'" + $v6[1, $colAA] }

# Rotate: new row 4 <- old row 5, new row 5 <- old row 6, new row 6 <- old row 4
$rng4.Value2 = $v5
$rng5.Value2 = $v6
$rng6.Value2 = $v4
